$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-fetch layout: Date moves to column A, MLBSO00 moves to column C ---
# (LNBSF00 in column B is unaffected by the reorder.)

$lastRow = 31

# Capture the existing column A (MLBSO00) and column C (Date) values before
# we start overwriting anything. Value2 is used so dates/numbers come back
# as plain numbers rather than formatted/boxed variants.
$oldA = @{}
$oldC = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $oldA[$r] = $ws.Cells.Item($r, 1).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
}

# Swap the header cells: A1 "MLBSO00" <-> C1 "Date".
$ws.Range("A1").Value = "Date"
$ws.Range("C1").Value = "MLBSO00"

# Swap the data columns for the existing rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $oldC[$r]
    $ws.Cells.Item($r, 3).Value = $oldA[$r]
}

# Column A is now the date column - apply the same date/time format used by
# the old Date column (except the very last row, handled below).
$ws.Range("A2:A$lastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Column C is now plain numeric MLBSO00 data - no special formatting.
$ws.Range("C2:C$lastRow").ClearFormats()

# --- Append the newest day's row (row 32), duplicating the latest reading ---
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = $oldC[$lastRow]
$ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($lastRow, 2).Value2
$ws.Cells.Item($newRow, 3).Value = $oldA[$lastRow]

# The "latest row" date-only format moves from the old last row to the new one.
$ws.Range("A$lastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A$newRow").NumberFormat = "YYYY-MM-DD"
